# Rename the auto-generated heading bookmarks (the "X<hash>" placeholder
# names) to human-readable slugs derived from each heading's text, per the
# commit "mlearning 2019 work in progress".
#
# The runtime doesn't support Bookmark.Delete / Bookmark.Name= as real
# mutations, so we rebuild each target paragraph's OOXML in place (same
# pStyle, same run text) with a freshly named bookmark wrapping it. This
# removes the old bookmark and adds the new one in a single atomic
# replace of that paragraph's content.

$d = $word.ActiveDocument

$renames = @(
    @{ Old = "X1138f39a1f8a0a41bf7a7cb18d9bac037b92eb3"; New = "introduction-to-online-blended-and-web-enhanced-learning" },
    @{ Old = "X265d2a0abd4b66ae60a117f8e09fac76057761b"; New = "evaluation-of-online-blended-and-web-enhanced-models" },
    @{ Old = "X3c799ec964ac7fdb78fe3becb6ce744659d0883"; New = "the-notion-of-anytime-and-anyplace-learning" },
    @{ Old = "X1ddcd1fef09bbd46b73fe04bb6f88914ebd3ab4"; New = "technology-for-online-and-blended-learning" },
    @{ Old = "Xd7a62706af5272d8043f84302b2271828477f61"; New = "technology-for-online-and-blended-learning---part-2" },
    @{ Old = "Xaa15f1943583ea0030132a53798722da47861e6"; New = "online-blended-and-web-enhanced-kearning-in-k-12-and-higher-education" },
    @{ Old = "X94cb49c95b3210bfcb2488e112e5763e8e51bc7"; New = "interactivity-communication-and-active-learning" }
)

$xmlTemplate = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkStart w:id="0" w:name="{0}"/><w:r><w:t xml:space="preserve">{1}</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

foreach ($item in $renames) {
    $oldName = $item.Old
    $newName = $item.New

    # Locate the bookmark by name (numeric index iteration -- Bookmarks.Item
    # by name string is unreliable in this host) to recover its Start/End
    # and the heading text it wraps.
    $target = $null
    for ($i = 1; $i -le $d.Bookmarks.Count; $i++) {
        $bm = $d.Bookmarks.Item($i)
        if ($bm.Name -eq $oldName) {
            $target = $bm
            break
        }
    }

    if ($null -eq $target) {
        Write-Output "WARNING: bookmark not found: $oldName"
        continue
    }

    $headingText = $d.Range($target.Start, $target.End).Text

    # Expand to the whole paragraph (so the stale bookmarkStart/bookmarkEnd
    # pair is fully replaced rather than left behind as an empty marker).
    # NOTE: deriving the paragraph via `<range>.Paragraphs.Item(1)` is
    # unreliable in this host (it can bind to the document's first
    # paragraph instead of the one actually covered by <range>), so instead
    # walk Document.Paragraphs and pick the one whose span contains the
    # bookmark.
    $foundPara = $null
    for ($j = 1; $j -le $d.Paragraphs.Count; $j++) {
        $cand = $d.Paragraphs.Item($j)
        $cr = $cand.Range
        if ($cr.Start -le $target.Start -and $target.End -le $cr.End) {
            $foundPara = $cand
            break
        }
    }
    $pRange = $foundPara.Range

    $xml = [string]::Format($xmlTemplate, $newName, $headingText)
    $pRange.InsertXML($xml)

    Write-Output "Renamed '$oldName' -> '$newName'"
}
